# Update profit_data.xlsx after running on 2026-01-16
# - adds two new "Combined Total Profit" columns (K, L) with header cells
# - pads existing rows 2-52 with (empty) cells in the new columns
# - appends a new data row (53) for 01/16/2026, including the combined totals

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells K1 / L1 -------------------------------------------
# Copy the formatting of the existing header cell (J1) onto the two new
# header cells so they pick up the same bold/centered/bordered style,
# then set their text.
$ws.Range("J1").Copy()
$ws.Range("K1:L1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("K1").Value = "Combined Total Profit(USD)"
$ws.Range("L1").Value = "Combined Total Profit(%)"

# --- Blank placeholder cells for the new columns on existing rows ------
# Touch the cells (without giving them a custom number format) so they
# exist in the sheet, matching the padded/empty cells added for rows 2-52.
$ws.Range("K2:L52").Style = "Normal"

# --- New row of data (row 53, 01/16/2026) -------------------------------
# Column A holds the date as plain text (matching the rest of the column),
# so format it as Text before assigning the value to avoid Excel turning
# it into a date serial number.
$ws.Range("A53").NumberFormat = "@"
$ws.Range("A53").Value = "01/16/2026"

$ws.Range("B53").Value = 12796.09
$ws.Range("C53").Value = 0.2266795759358105
$ws.Range("D53").Value = 0.7733204240641895
$ws.Range("E53").Value = -137.53
$ws.Range("F53").Value = -20.08
$ws.Range("G53").Value = -20846.28
$ws.Range("H53").Value = -67.81
$ws.Range("I53").Value = -252.24
$ws.Range("J53").Value = -8
$ws.Range("K53").Value = -21098.52
$ws.Range("L53").Value = -62.25
